# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (week of 2022-07-11) at the top of the date-descending
# price table, pushing the existing rows 698:732 down to 701:735.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 698-732 down by 3 rows, creating blank rows 698:700
$ws.Rows("698:700").Insert()

# Row 698 - Packham's Triumph / Especial
$ws.Cells.Item(698, 1).Value = 8
$ws.Cells.Item(698, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(698, 3).Value = "Coquimbo"
$ws.Cells.Item(698, 4).Value = 44753
$ws.Cells.Item(698, 5).Value = 4
$ws.Cells.Item(698, 6).Value = "Fruta"
$ws.Cells.Item(698, 7).Value = 100104
$ws.Cells.Item(698, 8).Value = "Frutos de pepita"
$ws.Cells.Item(698, 9).Value = 100104005
$ws.Cells.Item(698, 10).Value = "Pera"
$ws.Cells.Item(698, 11).Value = "Packham's Triumph"
$ws.Cells.Item(698, 12).Value = "Especial"
$ws.Cells.Item(698, 13).Value = 20
$ws.Cells.Item(698, 14).Value = 200000
$ws.Cells.Item(698, 15).Value = 210000
$ws.Cells.Item(698, 16).Value = 205000
$ws.Cells.Item(698, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(698, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(698, 19).Value = 456
$ws.Cells.Item(698, 20).Value = 450

# Row 699 - Packham's Triumph / Primera
$ws.Cells.Item(699, 1).Value = 8
$ws.Cells.Item(699, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(699, 3).Value = "Coquimbo"
$ws.Cells.Item(699, 4).Value = 44753
$ws.Cells.Item(699, 5).Value = 4
$ws.Cells.Item(699, 6).Value = "Fruta"
$ws.Cells.Item(699, 7).Value = 100104
$ws.Cells.Item(699, 8).Value = "Frutos de pepita"
$ws.Cells.Item(699, 9).Value = 100104005
$ws.Cells.Item(699, 10).Value = "Pera"
$ws.Cells.Item(699, 11).Value = "Packham's Triumph"
$ws.Cells.Item(699, 12).Value = "Primera"
$ws.Cells.Item(699, 13).Value = 20
$ws.Cells.Item(699, 14).Value = 180000
$ws.Cells.Item(699, 15).Value = 190000
$ws.Cells.Item(699, 16).Value = 185000
$ws.Cells.Item(699, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(699, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(699, 19).Value = 411
$ws.Cells.Item(699, 20).Value = 450

# Row 700 - Packham's Triumph / Segunda
$ws.Cells.Item(700, 1).Value = 8
$ws.Cells.Item(700, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(700, 3).Value = "Coquimbo"
$ws.Cells.Item(700, 4).Value = 44753
$ws.Cells.Item(700, 5).Value = 4
$ws.Cells.Item(700, 6).Value = "Fruta"
$ws.Cells.Item(700, 7).Value = 100104
$ws.Cells.Item(700, 8).Value = "Frutos de pepita"
$ws.Cells.Item(700, 9).Value = 100104005
$ws.Cells.Item(700, 10).Value = "Pera"
$ws.Cells.Item(700, 11).Value = "Packham's Triumph"
$ws.Cells.Item(700, 12).Value = "Segunda"
$ws.Cells.Item(700, 13).Value = 16
$ws.Cells.Item(700, 14).Value = 160000
$ws.Cells.Item(700, 15).Value = 170000
$ws.Cells.Item(700, 16).Value = 165000
$ws.Cells.Item(700, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(700, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(700, 19).Value = 367
$ws.Cells.Item(700, 20).Value = 450
